# chooseBlockD.xlsx — "delete visual_clutter and restructure stimulus folder"
#
# The visual_clutter (non-social "scrambled") stimulus set is replaced by a
# simpler Social1/Social2/Nonsocial1/Nonsocial2 PNG set, and the stimuli are
# now stored flat under stimuli/ instead of stimuli/social/ and
# stimuli/non-social/. Only the second (data) row of the sheet changes;
# the header row (posFile/stimFile/cs_plus_s/cs_minus_s/cs_plus_ns/cs_minus_ns)
# and the first two data cells (positions.xlsx / stimuli.xlsx) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "stimuli/Social2.png"
$ws.Range("D2").Value = "stimuli/Social1.png"
$ws.Range("E2").Value = "stimuli/Nonsocial2.png"
$ws.Range("F2").Value = "stimuli/Nonsocial1.png"

# The merged C:D / E:F column groups keep the same "best fit" width class as
# column A/B (16.53515625) and a wider class for the longer Nonsocial paths
# (19.84375) respectively.
$ws.Columns.Item(3).ColumnWidth = 15.65
$ws.Columns.Item(4).ColumnWidth = 15.65
$ws.Columns.Item(5).ColumnWidth = 19.0
$ws.Columns.Item(6).ColumnWidth = 19.0

# Selection moved off the original B1:B2 block onto F5.
$ws.Range("F5").Select() | Out-Null
